$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" (strikeouts) column values, replacing the old "Strike#" derived values.
$kValues = @{
    2  = 4
    3  = 4
    4  = 3
    5  = 2
    6  = 4
    7  = 6
    8  = 2
    9  = 7
    10 = 4
    11 = 3
    12 = 4
    13 = 2
    14 = 3
    15 = 3
    16 = 3
    17 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
